$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8296735286712646
$ws.Range("B1").Value = 1.698322415351868
$ws.Range("C1").Value = 6.595419883728027
$ws.Range("D1").Value = 2.133639812469482
$ws.Range("E1").Value = 0.950488269329071
